$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 130, shifting existing rows 130:212 down to 131:213
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new data record
$ws.Cells.Item(130, 1).Value = 5
$ws.Cells.Item(130, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(130, 3).Value = "Maule"
$ws.Cells.Item(130, 4).Value = 44777
$ws.Cells.Item(130, 5).Value = 7
$ws.Cells.Item(130, 6).Value = 100112017
$ws.Cells.Item(130, 7).Value = "Apio"
$ws.Cells.Item(130, 8).Value = "Americana (o)"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 600
$ws.Cells.Item(130, 11).Value = 9000
$ws.Cells.Item(130, 12).Value = 9000
$ws.Cells.Item(130, 13).Value = 9000
$ws.Cells.Item(130, 14).Value = "`$/docena de matas"
$ws.Cells.Item(130, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(130, 16).Value = 1500
$ws.Cells.Item(130, 17).Value = 6
$ws.Cells.Item(130, 18).Value = "Hortaliza"
